$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the previously empty "Status" cell for the bank-account-creation
# task as DONE; this also recalculates the dependent shared formula in D3
# and the total in D13.
$ws.Range("C3").Value = "DONE"

# Update the active selection to C4 (matches saved cursor position).
$ws.Range("C4").Select() | Out-Null
